$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts the whole sheet right by one).
$ws.Range("A1").EntireColumn.Insert()

# Copy the (now shifted) header formatting onto the new index column so the
# new A2:A8 cells pick up the same bold/bordered/centered style as row 1.
$ws.Range("B1").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)

# New leading index column: 0..6 for data rows 2..8.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6

# New row for the "KS" user (row 8). Name goes in column C; the rest of the
# row (besides the new index cell already set above) stays blank, matching
# the empty-string placeholders used elsewhere in the sheet for missing data.
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "KS"
$ws.Range("D8:V8").Value = ""

Write-Host "done"
